$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 45725
$ws.Range("B11").Value = "1 hours"
$ws.Range("C11").Value = "discuss plots and next steps"
$ws.Range("D11").Value = "Grant"

$ws.Range("A12").Select()
